# Regenerate merged AHB files
# - Rename header columns: *_old -> *_FV2410, *_new -> *_FV2504
# - Add a Table (ListObject) over the used range A1:U69
# - Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row's "_old"/"_new" suffixes.
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the used range into an Excel Table ("ListObject").
$rng = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes, $null)
$tbl.Name = "Table1"

# 3. Freeze the top (header) row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
